$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 310-311 for a new reporting date (45021 = 2023-04-05),
# pushing the existing rows 310-328 down to 312-330.
$ws.Rows("310:311").Insert()

# Row 310: Artic Sprite / Especial
$ws.Cells.Item(310, 1).Value = 2
$ws.Cells.Item(310, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(310, 3).Value = "Coquimbo"
$ws.Cells.Item(310, 4).Value = 45021
$ws.Cells.Item(310, 5).Value = 4
$ws.Cells.Item(310, 6).Value = "Fruta"
$ws.Cells.Item(310, 7).Value = 100103
$ws.Cells.Item(310, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(310, 9).Value = 100103006
$ws.Cells.Item(310, 10).Value = "Nectarín"
$ws.Cells.Item(310, 11).Value = "Artic Sprite"
$ws.Cells.Item(310, 12).Value = "Especial"
$ws.Cells.Item(310, 13).Value = 10
$ws.Cells.Item(310, 14).Value = 550000
$ws.Cells.Item(310, 15).Value = 560000
$ws.Cells.Item(310, 16).Value = 555000
$ws.Cells.Item(310, 17).Value = "$/bins (420 kilos)"
$ws.Cells.Item(310, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(310, 19).Value = 1321
$ws.Cells.Item(310, 20).Value = 420

# Row 311: Artic Sprite / Primera
$ws.Cells.Item(311, 1).Value = 2
$ws.Cells.Item(311, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(311, 3).Value = "Coquimbo"
$ws.Cells.Item(311, 4).Value = 45021
$ws.Cells.Item(311, 5).Value = 4
$ws.Cells.Item(311, 6).Value = "Fruta"
$ws.Cells.Item(311, 7).Value = 100103
$ws.Cells.Item(311, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(311, 9).Value = 100103006
$ws.Cells.Item(311, 10).Value = "Nectarín"
$ws.Cells.Item(311, 11).Value = "Artic Sprite"
$ws.Cells.Item(311, 12).Value = "Primera"
$ws.Cells.Item(311, 13).Value = 16
$ws.Cells.Item(311, 14).Value = 500000
$ws.Cells.Item(311, 15).Value = 510000
$ws.Cells.Item(311, 16).Value = 505000
$ws.Cells.Item(311, 17).Value = "$/bins (420 kilos)"
$ws.Cells.Item(311, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(311, 19).Value = 1202
$ws.Cells.Item(311, 20).Value = 420
